# Updated symbol list on Sat Jan 21 05:53:33 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto rows that
# moved since the last scrape. Values are written with a leading apostrophe so
# the numeric-looking/percentage-looking text is stored as literal text (as it
# already was in the workbook), not re-interpreted as a Number/Percentage by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.22"
$ws.Range("E2").Value = "'4.62%"

$ws.Range("E3").Value = "'12.39%"

$ws.Range("D4").Value = "'5.148"
$ws.Range("E4").Value = "'4.71%"

$ws.Range("D5").Value = "'0.07748"
$ws.Range("E5").Value = "'4.72%"

$ws.Range("D6").Value = "'2.355"
$ws.Range("E6").Value = "'3.90%"

$ws.Range("D7").Value = "'8.011"
$ws.Range("E7").Value = "'3.87%"

$ws.Range("D8").Value = "'3.950"
$ws.Range("E8").Value = "'5.34%"

$ws.Range("D9").Value = "'0.9295"
$ws.Range("E9").Value = "'1.93%"

$ws.Range("D10").Value = "'0.09960"
$ws.Range("E10").Value = "'11.93%"

$ws.Range("D11").Value = "'0.1796"
$ws.Range("E11").Value = "'6.58%"

$ws.Range("D12").Value = "'0.08613"
$ws.Range("E12").Value = "'4.55%"

$ws.Range("D13").Value = "'0.03322"
$ws.Range("E13").Value = "'6.42%"

$ws.Range("D14").Value = "'0.09889"
$ws.Range("E14").Value = "'-0.61%"

$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'-0.16%"

$ws.Range("D16").Value = "'0.005762"
$ws.Range("E16").Value = "'-1.65%"

$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'-0.79%"

$ws.Range("D18").Value = "'2.135"
$ws.Range("E18").Value = "'2.31%"

$ws.Range("D19").Value = "'0.3367"
$ws.Range("E19").Value = "'1.19%"

$ws.Range("E20").Value = "'2.77%"

$ws.Range("D21").Value = "'4.354"
$ws.Range("E21").Value = "'9.67%"

$ws.Range("E22").Value = "'5.14%"

$ws.Range("D23").Value = "'0.04579"
$ws.Range("E23").Value = "'0.47%"

$ws.Range("E24").Value = "'0.58%"

$ws.Range("D25").Value = "'0.004461"
$ws.Range("E25").Value = "'-2.61%"

$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'-0.01%"

$ws.Range("E27").Value = "'-0.21%"

$ws.Range("D39").Value = "'0.01790"
$ws.Range("E39").Value = "'12.41%"

$ws.Range("D40").Value = "'0.04790"
$ws.Range("E40").Value = "'7.06%"

$ws.Range("D41").Value = "'0.007788"
$ws.Range("E41").Value = "'6.67%"

$ws.Range("E42").Value = "'6.46%"

$ws.Range("D43").Value = "'0.007161"
$ws.Range("E43").Value = "'-25.23%"

$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'-5.91%"

$ws.Range("D45").Value = "'0.009180"
$ws.Range("E45").Value = "'3.32%"

$ws.Range("D46").Value = "'0.00006121"
$ws.Range("E46").Value = "'0.74%"

$ws.Range("E47").Value = "'-0.09%"

$ws.Range("E48").Value = "'39.20%"

$ws.Range("D49").Value = "'0.002000"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.09%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.09%"
